$d = $word.ActiveDocument

# The grade-band label uses an en dash (U+2013): "50 - 50%" (visually).
# NOTE: this interpreter coerces "+" to numeric addition whenever an
# operand looks number-like (even mixed with [char]), so build the
# needle via string interpolation instead of "+" concatenation.
$dash = [char]0x2013
$needle = "50 $dash 50%"

$rng = $d.Content
$found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target text '$needle' in document"
}

# $rng now spans exactly the 8 characters "50 - 50%", which live inside
# a single bold run: <w:r><w:rPr><w:b/></w:rPr><w:t>50 - 50%</w:t></w:r>,
# itself the sole content of its <w:p>. The target edit turns that one
# run into three bold runs - "50 - 5", "9", "%" - so the visible text
# becomes "50 - 59%" (fixing what was a typo'd "50 - 50%" in a grade
# band table that otherwise reads 70-100% / 60-69% / 50-59% / 40-49%).
#
# A plain Range.Text assignment gets re-coalesced back into a single
# run on save, so the run split is rebuilt explicitly via InsertXML,
# which replaces the matched range with literal OOXML. The surrounding
# paragraph's own identity attributes/<w:pPr> (paraId/textId/rsid*,
# BodyText style, centered, bold) are reproduced unchanged so only the
# runs themselves change, matching the target edit exactly.
$dash = [char]0x2013
$run1 = "<w:r><w:rPr><w:b/></w:rPr><w:t>50 $dash 5</w:t></w:r>"
$run2 = "<w:r><w:rPr><w:b/></w:rPr><w:t>9</w:t></w:r>"
$run3 = "<w:r><w:rPr><w:b/></w:rPr><w:t>%</w:t></w:r>"

$pOpen = '<w:p w14:paraId="2648481E" w14:textId="1FBAB7C8" w:rsidR="00580E83" w:rsidRPr="00352842" w:rsidRDefault="00580E83" w:rsidP="008D5336">'
$pPr = '<w:pPr><w:pStyle w:val="BodyText"/><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr>'
$newParaXml = "$pOpen$pPr$run1$run2$run3</w:p>"

$wrapped = "<?xml version=`"1.0`" standalone=`"yes`"?>" +
    "<?mso-application progid=`"Word.Document`"?>" +
    "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
    "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
    "<pkg:xmlData>" +
    "<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`" xmlns:w14=`"http://schemas.microsoft.com/office/word/2010/wordml`">" +
    "<w:body>$newParaXml</w:body>" +
    "</w:document>" +
    "</pkg:xmlData></pkg:part></pkg:package>"

$rng.InsertXML($wrapped)
